# Add team record (Wins/Losses/Ties) columns to the sheet, mirroring the
# existing header style used by the other header cells (e.g. AC1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing header cell onto the new header cells
# so they end up sharing the same cell style as the rest of row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122, $false) | Out-Null

# Set header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2 through 54) gets the same team record values.
for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 30).Value = 77
    $ws.Cells.Item($row, 31).Value = 85
    $ws.Cells.Item($row, 32).Value = 0
}
